# thermometer adjust gameflow, some fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New localization rows describing the "broken thermometer" mini-flow.
# Cells are written in the same order the strings were originally typed in
# (this matters for shared-string table ordering on save).
$ws.Cells.Item(62, 1).Value = "thermometer_broken"
$ws.Cells.Item(62, 2).Value = "Something is not right with this thermometer, should we proceed to use it?"
$ws.Cells.Item(62, 3).Value = 3

$ws.Cells.Item(63, 2).Value = "Yeah, it's fine."
$ws.Cells.Item(63, 1).Value = "thermometer_broken0"
$ws.Cells.Item(63, 3).Value = 1.5

$ws.Cells.Item(64, 1).Value = "thermometer_broken1"
$ws.Cells.Item(64, 2).Value = "No, it's broken."
$ws.Cells.Item(64, 3).Value = 1.5

$ws.Cells.Item(65, 1).Value = "thermometer_broken_correct"

$ws.Cells.Item(66, 1).Value = "thermometer_broken_wrong"
$ws.Cells.Item(66, 2).Value = "Wrong! The arrow wasn't moving at all when the thermometer was dipped into the glass of ice."

$ws.Cells.Item(65, 2).Value = "That's right, the arrow wasn't moving at all when the thermometer was dipped into the glass of ice. Good call!"

# Leave the view roughly where the author left it after the edit.
$ws.Range("B68").Select()
